# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a brand-new data row before the current row 244 -------------
# This pushes the existing rows 244..329 down to 245..330.
$ws.Rows.Item(244).Insert()

$newRow = 244
$ws.Cells.Item($newRow, 1).Value = 5
$ws.Cells.Item($newRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($newRow, 3).Value = "Maule"
$ws.Cells.Item($newRow, 4).Value = 45120
$ws.Cells.Item($newRow, 5).Value = 7
$ws.Cells.Item($newRow, 6).Value = 100112017
$ws.Cells.Item($newRow, 7).Value = "Apio"
$ws.Cells.Item($newRow, 8).Value = "Americana (o)"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 600
$ws.Cells.Item($newRow, 11).Value = 6000
$ws.Cells.Item($newRow, 12).Value = 6000
$ws.Cells.Item($newRow, 13).Value = 6000
$ws.Cells.Item($newRow, 14).Value = "`$/docena de matas"
$ws.Cells.Item($newRow, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($newRow, 16).Value = 1000
$ws.Cells.Item($newRow, 17).Value = 6
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Make sure the date-formatted column keeps its expected number format
# (same as neighbouring rows in column D).
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat

# --- 2) Append one more row at the very end, duplicating what is now the
#        last existing row (originally row 329, now shifted to row 330). ---
$lastRow = 331
$ws.Cells.Item($lastRow, 1).Value = 5
$ws.Cells.Item($lastRow, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item($lastRow, 3).Value = "Maule"
$ws.Cells.Item($lastRow, 4).Value = 45072
$ws.Cells.Item($lastRow, 5).Value = 7
$ws.Cells.Item($lastRow, 6).Value = 100112017
$ws.Cells.Item($lastRow, 7).Value = "Apio"
$ws.Cells.Item($lastRow, 8).Value = "Americana (o)"
$ws.Cells.Item($lastRow, 9).Value = "Primera"
$ws.Cells.Item($lastRow, 10).Value = 700
$ws.Cells.Item($lastRow, 11).Value = 6000
$ws.Cells.Item($lastRow, 12).Value = 6000
$ws.Cells.Item($lastRow, 13).Value = 6000
$ws.Cells.Item($lastRow, 14).Value = "`$/docena de matas"
$ws.Cells.Item($lastRow, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($lastRow, 16).Value = 1000
$ws.Cells.Item($lastRow, 17).Value = 6
$ws.Cells.Item($lastRow, 18).Value = "Hortaliza"

$ws.Cells.Item($lastRow, 4).NumberFormat = $ws.Cells.Item($lastRow - 1, 4).NumberFormat
